$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 1468.566
$ws.Range("I15").Value = 1468.566
$ws.Range("K15").Value = 4405.698
$ws.Range("M15").Value = -4236.698
# row 62
$ws.Range("H62").Value = 21498.63
$ws.Range("I62").Value = 21643.428
$ws.Range("J62").Value = 21414.166
$ws.Range("K62").Value = 21643.428
$ws.Range("L62").Value = 21414.166
$ws.Range("M62").Value = -21019.428
$ws.Range("N62").Value = -22662.166
# row 65
$ws.Range("H65").Value = 21498.63
$ws.Range("I65").Value = 21643.428
$ws.Range("J65").Value = 21414.166
$ws.Range("K65").Value = 108217.14
$ws.Range("L65").Value = 107070.83
$ws.Range("M65").Value = -105097.14
$ws.Range("N65").Value = -113310.83
# row 76
$ws.Range("H76").Value = 4237.6313
$ws.Range("I76").Value = 2893.9285
$ws.Range("K76").Value = 2893.9285
$ws.Range("M76").Value = -2578.9285
# row 79
$ws.Range("H79").Value = 4237.6313
$ws.Range("I79").Value = 2893.9285
$ws.Range("K79").Value = 2893.9285
$ws.Range("M79").Value = -1801.9285
# row 86
$ws.Range("H86").Value = 5158.1577
$ws.Range("J86").Value = 5200.6665
$ws.Range("L86").Value = 5200.6665
$ws.Range("N86").Value = -7446.6665
# row 89
$ws.Range("H89").Value = 5158.1577
$ws.Range("J89").Value = 5200.6665
$ws.Range("L89").Value = 26003.3325
$ws.Range("N89").Value = -37235.3325
# row 107
$ws.Range("H107").Value = 356.4
$ws.Range("J107").Value = 150
$ws.Range("L107").Value = 150
$ws.Range("N107").Value = -3990
# row 112
$ws.Range("H112").Value = 4967.5
$ws.Range("J112").Value = 5269.6826
$ws.Range("L112").Value = 15809.0478
$ws.Range("N112").Value = -18025.0478
# row 136
$ws.Range("H136").Value = 73779.664
$ws.Range("J136").Value = 73779.664
$ws.Range("L136").Value = 73779.664
$ws.Range("N136").Value = -83979.664
# row 138
$ws.Range("H138").Value = 3128.2454
$ws.Range("I138").Value = 1406.6072
$ws.Range("J138").Value = 5056.48
$ws.Range("K138").Value = 4219.821599999999
$ws.Range("L138").Value = 15169.44
$ws.Range("M138").Value = 920.1784000000007
$ws.Range("N138").Value = -25449.44

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 4
$ws.Range("H4").Value = 255.42857
$ws.Range("I4").Value = 255.42857
$ws.Range("K4").Value = 255.42857
$ws.Range("M4").Value = -139.42857
# row 32
$ws.Range("H32").Value = 34024.91
$ws.Range("I32").Value = 34024.91
$ws.Range("K32").Value = 34024.91
$ws.Range("M32").Value = -33737.91
# row 102
$ws.Range("H102").Value = 18287.523
$ws.Range("I102").Value = 20947.445
$ws.Range("K102").Value = 20947.445
$ws.Range("M102").Value = -19325.445
# row 132
$ws.Range("H132").Value = 1252205.1
$ws.Range("I132").Value = 1335418.8
$ws.Range("K132").Value = 4006256.4
$ws.Range("M132").Value = -4003726.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Range("H94").Value = 1862.5
$ws.Range("I94").Value = 1862.5
$ws.Range("K94").Value = 1862.5
$ws.Range("M94").Value = -1411.5
# row 95
$ws.Range("H95").Value = 67853.336
$ws.Range("J95").Value = 67853.336
$ws.Range("L95").Value = 67853.336
$ws.Range("N95").Value = -73345.336
# row 96
$ws.Range("H96").Value = 19928
$ws.Range("I96").Value = 19928
$ws.Range("K96").Value = 19928
$ws.Range("M96").Value = -17182
# row 97
$ws.Range("H97").Value = 9214
$ws.Range("J97").Value = 10000
$ws.Range("L97").Value = 10000
$ws.Range("N97").Value = -11982
# row 132
$ws.Range("H132").Value = 99989.5
$ws.Range("J132").Value = 99989.5
$ws.Range("L132").Value = 99989.5
$ws.Range("N132").Value = -110109.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 62
$ws.Range("H62").Value = 6461.125
$ws.Range("J62").Value = 8897.799999999999
$ws.Range("L62").Value = 8897.799999999999
$ws.Range("N62").Value = -10145.8
# row 65
$ws.Range("H65").Value = 6461.125
$ws.Range("J65").Value = 8897.799999999999
$ws.Range("L65").Value = 44489
$ws.Range("N65").Value = -50729
# row 93
$ws.Range("H93").Value = 39997
$ws.Range("I93").Value = 39997
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 39997
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -38125
$ws.Range("N93").ClearContents()
# row 95
$ws.Range("H95").Value = 15832.667
$ws.Range("J95").Value = 15832.667
$ws.Range("L95").Value = 15832.667
$ws.Range("N95").Value = -21324.667
# row 96
$ws.Range("H96").Value = 4923.485
$ws.Range("J96").Value = 5119.8276
$ws.Range("L96").Value = 5119.8276
$ws.Range("N96").Value = -10611.8276
# row 99
$ws.Range("H99").Value = 3422.4285
$ws.Range("I99").Value = 2420.6667
$ws.Range("K99").Value = 2420.6667
$ws.Range("M99").Value = -922.6667000000002
# row 126
$ws.Range("H126").Value = 3422.4285
$ws.Range("I126").Value = 2420.6667
$ws.Range("K126").Value = 7262.000100000001
$ws.Range("M126").Value = -4792.000100000001
# row 132
$ws.Range("H132").Value = 67011664
$ws.Range("I132").Value = 111114140
$ws.Range("K132").Value = 333342420
$ws.Range("M132").Value = -333339890

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 36
$ws.Range("H36").Value = 833.3333
$ws.Range("I36").Value = 725
$ws.Range("K36").Value = 2175
$ws.Range("M36").Value = -2006
# row 56
$ws.Range("H56").Value = 8426
$ws.Range("I56").Value = 8426
$ws.Range("K56").Value = 8426
$ws.Range("M56").Value = -7896
# row 57
$ws.Range("H57").Value = 1998
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
# row 139
$ws.Range("H139").Value = 3408.6365
$ws.Range("I139").Value = 2436.875
$ws.Range("K139").Value = 7310.625
$ws.Range("M139").Value = -2170.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 4624.75
$ws.Range("I70").Value = 3999.6667
$ws.Range("J70").Value = 6500
$ws.Range("K70").Value = 3999.6667
$ws.Range("L70").Value = 6500
$ws.Range("M70").Value = -3729.6667
$ws.Range("N70").Value = -7040
# row 73
$ws.Range("H73").Value = 4624.75
$ws.Range("I73").Value = 3999.6667
$ws.Range("J73").Value = 6500
$ws.Range("K73").Value = 3999.6667
$ws.Range("L73").Value = 6500
$ws.Range("M73").Value = -3063.6667
$ws.Range("N73").Value = -8372
# row 102
$ws.Range("H102").Value = 2809.1555
$ws.Range("I102").Value = 1829.9656
$ws.Range("K102").Value = 1829.9656
$ws.Range("M102").Value = -207.9656
# row 113
$ws.Range("H113").Value = 6565.2173
$ws.Range("I113").Value = 5833.3335
$ws.Range("J113").Value = 7363.636
$ws.Range("K113").Value = 5833.3335
$ws.Range("L113").Value = 7363.636
$ws.Range("M113").Value = -3663.3335
$ws.Range("N113").Value = -11703.636
# row 132
$ws.Range("H132").Value = 112453540
$ws.Range("I132").Value = 112453540
$ws.Range("K132").Value = 337360620
$ws.Range("M132").Value = -337358090
# row 137
$ws.Range("H137").Value = 49500
$ws.Range("J137").Value = 49500
$ws.Range("L137").Value = 49500
$ws.Range("N137").Value = -59700

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 621.625
$ws.Range("I16").Value = 396.66666
$ws.Range("K16").Value = 396.66666
$ws.Range("M16").Value = -226.66666
# row 46
$ws.Range("H46").Value = 804.4
$ws.Range("I46").Value = 442.33334
$ws.Range("K46").Value = 442.33334
$ws.Range("M46").Value = -254.33334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 110
$ws.Range("H110").Value = 148527.5
$ws.Range("J110").Value = 148527.5
$ws.Range("L110").Value = 148527.5
$ws.Range("N110").Value = -156707.5

Write-Host "All updates applied."